$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "paul"
$ws.Range("B5").Value = "paul@gmail.com"
$ws.Range("C5").Value = "test"
$ws.Range("D5").Value = "test"
